$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.492.34"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.914.88"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.09"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2888"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.65"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.36"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("D12").Value = "1.911.88"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07563"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.242"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6676"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "304.31"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "30.488.65"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007563"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "2.163.01"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.471"
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9984"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.412"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.477"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.98"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("E27").Value = "  -6.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.099"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1077"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.163"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7350"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02050"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.726"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.674"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.74"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.018"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4423"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8642"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.904"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9988"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "69.07"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.00"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.284"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.304"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2520"
$ws.Range("E51").Value = "  +0.69%  "